# Create button 'Save and Load'. Formatting datetime in filters table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new timesheet entry in row 29 (previously a blank placeholder row)
$ws.Range("A29").Value = "Работа по созданию функционала загрузки данных (Отложенная загрузка данных с Celery)"
$ws.Range("B29").Value = 3

# Copy the date-formatted style from the row above so the new date cell
# keeps the same number format (m/d/yyyy) instead of creating a new one
$ws.Range("C28").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = "4/24/2019"

# Leave the selection where the user ended up after entering the new row
$ws.Range("C30").Select()
